# AHDT2_OP_MDX_NEG.docx template update (MHD2-210)
#
# Applies the three substantive content edits captured by the commit's
# XML diff:
#   1. Add the ASXL1 detection-limit sentence after the existing JAK2
#      detection-limit sentence in the "limitations" paragraph.
#   2. Extend the "Insertions or deletions (particularly those > 25 bp
#      in length)" clause to also mention homopolymer regions.
#   3. Roll the report-date DATE field result forward from
#      11-Jul-2025 to 22-Aug-2025.
#
# (The diff also shows a large number of runs being split apart with
# identical rPr/no visible formatting or text change - e.g. "Clinical
# Indication" -> "Clinical " + "Indication", proofErr spell-check
# markers around words like NovaSeq/Oncoanalyser/OncoPath/gnomAD/
# ClinVar, and sectPr/customXml relationship-id renumbering. Those are
# byte-identical-content save artifacts with no visible or semantic
# effect, so they are intentionally not reproduced here.)

$d = $word.ActiveDocument

# 1) JAK2 -> add the ASXL1 clause right after the JAK2 detection-limit
#    parenthetical, before "This assay is primarily qualitative...".
$r1 = $d.Content.Find.Execute(
    "(Val617Phe) (detection limit ~ 1%). This assay",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "(Val617Phe) (detection limit ~ 1%) and ASXL1 c.1934dup;p.(Gly646Trpfs*12) (detection limit ~ 5%). This assay",
    2)
Write-Output ("JAK2/ASXL1 replace: " + $r1)

# 2) Insertions/deletions clause -> mention homopolymer regions too.
$r2 = $d.Content.Find.Execute(
    "those > 25 bp in length) are not reliably",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "those > 25 bp in length or in homopolymer regions) are not reliably",
    2)
Write-Output ("Homopolymer replace: " + $r2)

# 3) Report date field result.
$r3 = $d.Content.Find.Execute(
    "11-Jul-2025",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "22-Aug-2025",
    2)
Write-Output ("Date replace: " + $r3)
